# Auto-generated Word COM-interop PowerShell script
# Applies the diff: date text change, insert 6 new rows after the table
# header (before the current row 1), update text in the next 14 rows,
# and delete the last 6 rows of the table.

$d = $word.ActiveDocument

# 1) Update the date line
$null = $d.Content.Find.Execute("2023-09-25 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-26 Tuesday", 2)

$t = $d.Tables.Item(1)

# 2) New rows to insert before the current first row
$newRowsData = @(
    @("7+78=","1+38=","77-50=","92-79=","84-55="),
    @("14+43=","2+67=","17+27=","26-13=","98-48="),
    @("9+36=","64-43=","84+15=","82-36=","5+82="),
    @("88-75=","51+46=","62-28=","31-24=","76-59="),
    @("52+46=","40+40=","53-7=","68-32=","15+39="),
    @("8+76=","52-35=","49-46=","20-5=","48-2=")
)

$firstRow = $t.Rows.Item(1)
for ($i = $newRowsData.Length - 1; $i -ge 0; $i--) {
    $rowData = $newRowsData[$i]
    $newRow = $t.Rows.Add($firstRow)
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $newRow.Cells.Item($c + 1).Range.Text = $rowData[$c]
    }
}

# 3) Update text for the rows that follow the newly inserted ones
#    (these were rows 1-14 before the insert, now rows 7-20)
$updatedRowsData = @(
    @("98-71=","50-31=","16-2=","47+45=","53-3="),
    @("42+8=","36+63=","58-25=","6+22=","0+62="),
    @("32-25=","40-0=","88-49=","87-5=","25-16="),
    @("44+54=","77-1=","60+11=","84-80=","80-56="),
    @("25+23=","84-37=","30+66=","3+39=","73-12="),
    @("48+15=","78+0=","50+6=","43+43=","46+38="),
    @("16+30=","35-22=","7+2=","36+62=","15+59="),
    @("7+26=","84-28=","7+21=","16+15=","86-85="),
    @("64-28=","96-58=","16+32=","75-34=","48-40="),
    @("87-33=","15+41=","25+19=","32+25=","71-38="),
    @("95-49=","66-18=","43+48=","60+36=","1+63="),
    @("60-7=","78-6=","68-8=","84-53=","68+6="),
    @("46+4=","14+53=","97-13=","70-1=","44-13="),
    @("11+83=","61-26=","47+39=","8+56=","75-57=")
)

for ($i = 0; $i -lt $updatedRowsData.Length; $i++) {
    $rowIndex = 7 + $i
    $rowData = $updatedRowsData[$i]
    $row = $t.Rows.Item($rowIndex)
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $row.Cells.Item($c + 1).Range.Text = $rowData[$c]
    }
}

# 4) Delete the trailing 6 rows (were rows 15-20 before the insert)
for ($i = 0; $i -lt 6; $i++) {
    $t.Rows.Item($t.Rows.Count).Delete()
}
